$d = $word.ActiveDocument

$pairs = @(
    @("13×21=", "42×29="),
    @("72×70=", "41×14="),
    @("86×35=", "20×90="),
    @("24×39=", "11×14="),
    @("24×95=", "47×62="),
    @("88×57=", "78×44="),
    @("79×78=", "49×39="),
    @("81×97=", "18×76="),
    @("98×76=", "30×95="),
    @("33×72=", "60×56="),
    @("88×62=", "27×66="),
    @("16×30=", "40×43="),
    @("39×75=", "40×25="),
    @("62×27=", "84×45="),
    @("90×31=", "81×65="),
    @("39×68=", "55×74="),
    @("41×29=", "79×96="),
    @("11×56=", "13×50="),
    @("28×27=", "71×36="),
    @("22×83=", "41×15="),
    @("62×86=", "22×32="),
    @("11×15=", "63×79="),
    @("70×20=", "79×24="),
    @("91×15=", "42×45="),
    @("83×40=", "48×82=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
